$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "adminNo"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "gender"
$ws.Range("D1").Value = "citizenshipStatus"
$ws.Range("E1").Value = "course"
$ws.Range("F1").Value = "stage"
$ws.Range("G1").Value = "pemGroup"

# --- Row 2 ---
$ws.Range("A2").Value = "212345A"
$ws.Range("B2").Value = "XLSX Uno"
$ws.Range("C2").Value = "Male"
$ws.Range("D2").Value = "Singapore citizen"
$ws.Range("E2").Value = "C02"
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "MI2002"

# --- Rows 3 & 4, filled column-by-column (matches the admin/name entry order
#     that produced the shared-strings layout: adminNo column down, then the
#     name column down, before the remaining repeated columns) ---
$ws.Range("A3").Value = "212346A"
$ws.Range("A4").Value = "212347A"

$ws.Range("B3").Value = "XLSX Dos"
$ws.Range("B4").Value = "XLSX Tres"

$ws.Range("C3").Value = "Male"
$ws.Range("C4").Value = "Male"

$ws.Range("D3").Value = "Singapore citizen"
$ws.Range("D4").Value = "Singapore citizen"

$ws.Range("E3").Value = "C02"
$ws.Range("E4").Value = "C02"

$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1

$ws.Range("G3").Value = "MI2002"
$ws.Range("G4").Value = "MI2002"

# --- Column widths (closest representable values - engine quantizes to 1/6 char) ---
$ws.Columns.Item(4).ColumnWidth = 17
$ws.Columns.Item(7).ColumnWidth = 13

# --- Leave the cursor where data entry ended, like Excel would after Tab/Enter ---
$ws.Range("F5").Select()
